# Apply updated TPM-derived NATMI statistics (Adam12-Itgb1, YoungD4) to rows 2-21,
# then drop the now-removed "Resolving-Mac" sending-cluster block (old rows 22-26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns E:T, one row per data row (sheet rows 2..21)
$newRows = @()
$newRows += ,@(3, 1, 9.666904666666666, 29.000714, 0.08665276800447426, 0.08665276800447427, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 1541.713117994574, 13875.41806195117, 0.02585099103532437, 0.02585099103532437)
$newRows += ,@(3, 1, 9.666904666666666, 29.000714, 0.08665276800447426, 0.08665276800447427, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 1668.106539922286, 15012.95885930057, 0.0279703186709531, 0.0279703186709531)
$newRows += ,@(3, 1, 9.666904666666666, 29.000714, 0.08665276800447426, 0.08665276800447427, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 719.0988590517263, 6471.889731465538, 0.0120576376641582, 0.01205763766415821)
$newRows += ,@(3, 1, 9.666904666666666, 29.000714, 0.08665276800447426, 0.08665276800447427, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 564.6885079004335, 5082.196571103902, 0.00946853041924767, 0.009468530419247672)
$newRows += ,@(3, 1, 9.666904666666666, 29.000714, 0.08665276800447426, 0.08665276800447427, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 674.2300209327412, 6068.070188394672, 0.01130529021479092, 0.01130529021479092)
$newRows += ,@(3, 1, 98.09641766666668, 294.289253, 0.8793224320759494, 0.8793224320759495, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 15644.77349884986, 140802.9614896488, 0.2623269496432158, 0.2623269496432159)
$newRows += ,@(3, 1, 98.09641766666668, 294.289253, 0.8793224320759494, 0.8793224320759495, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 16927.37039364425, 152346.3335427982, 0.2838331562404547, 0.2838331562404547)
$newRows += ,@(3, 1, 98.09641766666668, 294.289253, 0.8793224320759494, 0.8793224320759495, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 7297.167444342401, 65674.50699908161, 0.1223567523589172, 0.1223567523589172)
$newRows += ,@(3, 1, 98.09641766666668, 294.289253, 0.8793224320759494, 0.8793224320759495, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 5730.264405479921, 51572.37964931929, 0.09608338415696159, 0.09608338415696162)
$newRows += ,@(3, 1, 98.09641766666668, 294.289253, 0.8793224320759494, 0.8793224320759495, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 6841.853935405551, 61576.68541864995, 0.1147221896764, 0.1147221896764)
$newRows += ,@(2, 0.6666666666666666, 0.05829266666666667, 0.174878, 0.0005225272302980695, 0.0005225272302980696, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 9.296726509859557, 83.67053858873601, 0.0001558847692603518, 0.0001558847692603518)
$newRows += ,@(2, 0.6666666666666666, 0.05829266666666667, 0.174878, 0.0005225272302980695, 0.0005225272302980696, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 10.058894946122, 90.530054515098, 0.0001686645848974248, 0.0001686645848974248)
$newRows += ,@(2, 0.6666666666666666, 0.05829266666666667, 0.174878, 0.0005225272302980695, 0.0005225272302980696, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 4.33625773052511, 39.026319574726, 0.000072709091211777, 0.00007270909121177702)
$newRows += ,@(2, 0.6666666666666666, 0.05829266666666667, 0.174878, 0.0005225272302980695, 0.0005225272302980696, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 3.405143641794889, 30.646292776154, 0.00005709644468261002, 0.00005709644468261004)
$newRows += ,@(2, 0.6666666666666666, 0.05829266666666667, 0.174878, 0.0005225272302980695, 0.0005225272302980696, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 4.065692920549334, 36.591236284944, 0.00006817234024590587, 0.0000681723402459059)
$newRows += ,@(3, 1, 3.737483333333333, 11.21245, 0.03350227268927818, 0.03350227268927819, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 596.0674364727112, 5364.6069282544, 0.009994683042425186, 0.009994683042425189)
$newRows += ,@(3, 1, 3.737483333333333, 11.21245, 0.03350227268927818, 0.03350227268927819, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 644.9345065625499, 5804.41055906295, 0.01081407166672269, 0.01081407166672269)
$newRows += ,@(3, 1, 3.737483333333333, 11.21245, 0.03350227268927818, 0.03350227268927819, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 278.0228101340722, 2502.20529120665, 0.004661804513760959, 0.00466180451376096)
$newRows += ,@(3, 1, 3.737483333333333, 11.21245, 0.03350227268927818, 0.03350227268927819, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 218.3236474939278, 1964.91282744535, 0.003660786555092869, 0.003660786555092871)
$newRows += ,@(3, 1, 3.737483333333333, 11.21245, 0.03350227268927818, 0.03350227268927819, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 260.6753198630667, 2346.0778787676, 0.004370926911276475, 0.004370926911276477)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $sheetRow = $i + 2
    $vals = $newRows[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($sheetRow, $j + 5).Value = $vals[$j]
    }
}

# Remove the obsolete "Resolving-Mac" sending-cluster rows (previously rows 22-26)
$ws.Rows("22:26").Delete()
